$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Text",
    "Министерство науки и высшего образования Российской Федерации Федеральное государственное бюджетное образовательное учреждение высшего образования «Московский государственный технический университет имени Н.",
    "Э.",
    "Баумана (национальный исследовательский университет)» (МГТУ им.",
    "Н.",
    "Баумана) ФАКУЛЬТЕТ «Информатика исистемы управления» КАФЕДРА «Программное обеспечение ЭВМ иинформационные технологии» РАСЧЕТНО-ПОЯСНИТЕЛЬНАЯ ЗАПИСКА К НАУЧНО-ИССЛЕДОВАТЕЛЬСКОЙ РАБОТЕ НА ТЕМУ: «Вейвлет-шум» Студент ИУ7-53Б(BIS_19) (Группа) (Подпись, дата)Авдейкина В.",
    "П.",
    "(И.",
    "О.",
    "Фамилия) Руководитель НИР (Подпись, дата)Кострицкий А.",
    "С.",
    "Фамилия) 2023 г.InСОДЕРЖАНИЕ ВВЕДЕНИЕ 3 1 Обзор вейвлет-шума 4 1.1 Математические основы .",
    ".",
    "4 1.2 Вейвлет-шум .",
    "5 1.2.1 Получение R↓.",
    "6 1.2.2 Получение R↓.",
    "6 1.2.3 Получение R↓↑.",
    "7 1.2.4 Получение N.",
    "7 СПИСОК ИСПОЛЬЗОВАННЫХ ИСТОЧНИКОВ 9 2InВВЕДЕНИЕ За последние 15 лет генерация поверхностей с текстурами (стилями заполнения, имитирующими сложную рельефную объемную плоскость, выполненную из какого-то материала [1]) не перестала быть актуальной задачей компьютерной графики и упоминается в [2—6].",
    "Таким образом, возникает необходимость нахождения методов, которые обеспечивают решение этой задачи.",
    "Существует ряд решений описанной задачи, одним из которых является использование метода вейвлет-шума (вейвлет-шум) [7; 8], который будет рассмотрен далее.",
    "3In1 Обзор вейвлет-шума В данном разделе будут рассмотрены математические основы вейвлетшума и его описание.",
    "1.1 Математические основы Суть вейвлет-шума заключается в использовании вейвлетов.",
    "Вейвлет — функция независимой переменной, имеющая вид короткой волны (всплеска) [9—11].",
    "Пример вейвлета можно наблюдать на рисунке 1.",
    "[рисунок 1 — пример графика вейвлет функции] Центр вейвлета — значение независимой переменной, через которое проходит вертикальная ось симметрии вейвлета.",
    "Пространство вейвлетов W— набор функций, которые могут быть представлены линейной комбинацией сдвигов и масштабирования одного вейвлета.",
    "[12; 13] Масштабирующая функция ξ— функция, для которой справедливо следующее равенство: ξ(x) =X kpkξ(2x−k), (1.1) при условии, что существуют такие числовые коэффициенты pk.",
    "[7; 12] Согласно [7], если ϕ(x)— функция из W, центр которой x= 0, то любая функция F(x)может быть представлена следующим образом: F(x) =X ifiϕ(x−i), (1.2) гдеϕ(x−i)— функция из Wс центром x=i,fi— некоторые числовые коэффициенты.",
    "Пространство S0— пространство, состоящее из функций вида 1.2.",
    "Функции G(x)вида G(x) =X igiϕ(2x−i) (1.3) аналогично составляют пространство S1.",
    "[7] 4InЕсли рассматриваемая в 1.3, 1.2 функция ϕявляется масштабирующей, тоS1расширяет S0, то есть включает все функции S0.",
    "[7; 12] Вейвлет-анализ — процесс определения того, расширяет ли пространство S1пространство S0.",
    "[7].",
    "Тогда функция F(x)из 1.2 представляется с помощью коэффициентов f↑ iв пространстве S1: F(x) =X if↑ iϕ(2x−i) (1.4) Коэффициенты f↑ iвыражаются с использованием формулы 1.1: f↑ i=X kpi−2kfk (1.5) Согласно [7], коэффициенты g↓ i, необходимые для представления функцииG(x)из 1.3 в пространстве S0, вычисляются с помощью вейвлет-анализа и являются равными g↓ i=X kak−2igk (1.6) а функция G(x)представляется как G(x) =G↓(x) +D(x) (1.7) гдеD(x)— функция из S1, которая не может быть выражена в S0через какие-либо коэффициенты.",
    "1.2 Вейвлет-шум Изображение Xпредставляется набором числовых коэффициентов (.",
    ", x i, .",
    "): X= (.",
    ") (1.8) Входные данные вейвлет-шума: функция B(x)изW, изображение R.",
    "Выходные данные вейвлет-шума: изображение N.",
    "Авторы [7] представляют алгоритм вейвлет-шума тремя семантическими 5Inшагами: 1) получение R(x); 2) получение R↓(x); 3) получение R↓↑(x); 4) получение N(x).",
    "После выполнения алгоритма для получения результата значение функцииN(x)вычисляется при заданном При заданном вычисляется значение функции N(x)и получается изображение N= (.",
    ", n i, .",
    ").",
    "Пример визуализации получения результата алгоритма можно наблюдать на рис.",
    "2.",
    "[рисунок 2 — пример визуализации получения из R N[7]] Рассмотрим представленные шаги более подробно, основываясь на информации из [7].",
    "1.2.1 Получение R↓ Согласно формуле 1.8 изображение Rпредставляется следующим образом: R= (.",
    ", r i, .",
    ") (1.9) Аналогично формуле 1.3 выражается R(x): R=X iriB(2x−i) (1.10) 1.2.2 Получение R↓ Согласно формуле 1.7: R(x) =R↓(x) +N(x) (1.11) N(x)выражается из формулы 1.11 как N(x) =R(x)−R↓(x) (1.12) 6InС помощью формулы 1.2 получаем R↓(x) =X ir↓ iB(x−i) (1.13) гдеr↓ iвыражаются с помощью формулы 1.6: r↓ i=X kak−2irk (1.14) 1.2.3 Получение R↓↑ C помощью формул 1.10 и 1.13 выражение 1.12 преобразовывается как N(x) =X iriB(2x−i)−X ir↓ iB(x−i) (1.15) С учетом формулы 1.4 и 1.13 получается выражение R↓↑(x) =X ir↓↑ iB(2x−i) (1.16) где коэффициенты r↓↑ iвыражаются с учетом формулы 1.5: r↓↑ i=X kpi−2kr↓ k (1.17) 1.2.4 Получение N С использованием формулы 1.16 выражение 1.15 записывается как N(x) =X iriB(2x−i)−X ir↓↑ iB(2x−i) =X in′ iB(2x−i)(1.18) гдеn′ i=ri−r↓↑ i.",
    "Для получения изображения N= (.",
    ")значение функции N(x) вычисляется при задаваемых значениях x.",
    "7InСПИСОК ИСПОЛЬЗОВАННЫХ ИСТОЧНИКОВ 1.Порев В.",
    "Компьютерная графика.",
    "— СПб : БХВ-Петербург, 2002.",
    "— С.",
    "432.",
    "2.Pai H.",
    "-Y.Texture designs and workflows for physically based rendering using procedural texture generation //2019 IEEE Eurasia Conference on IOT, Communication and Engineering (ECICE).",
    "— 2019.",
    "195—198.",
    "3.Weakly-Supervised Photo-realistic Texture Generation for 3D Face Reconstruction / X.",
    "Yin [и др.] //2023 IEEE 17th International Conference on Automatic Face and Gesture Recognition (FG).",
    "— 2023.",
    "1—8.",
    "4.Digital preservation of Brazilian indigenous artworks: Generating high quality textures for 3D models [Электронный ресурс] / B.",
    "Trinch˜ ao Andrade [и др.] // Journal of Cultural Heritage.",
    "— 2012.",
    "— Т.",
    "13, № 1.",
    "28—39.",
    "— DOI: https://doi.org/10.1016/j.culher.2011.05.002 .",
    "— URL: https: //www.sciencedirect.com/science/article/pii/S1296207411000513 (дата обращения: 24.11.2023).",
    "5.Quintana G.",
    ",Ciurana J.",
    ",Ribatallada J.",
    "Surface Roughness Generation and Material Removal Rate in Ball End Milling Operations //Materials and Manufacturing Processes.",
    "— 2010.",
    "25, № 6.",
    "386—398.",
    "— (дата обращения: 24.11.2023).",
    "6.A Papier-Mˆ ach´ e Approach to Learning 3D Surface Generation / T.",
    "Groueix [и др.]//Proceedings of the IEEE Conference on Computer Vision and Pattern Recognition (CVPR).",
    "— 2018.",
    "7.Cook R.,De Rose T.",
    "Wavelet Noise [Электронный ресурс].",
    "— URL: https: //graphics.pixar.com/library/WaveletNoise/paper.pdf (дата обращения: 11.11.2023).",
    "8.Procedural Noise/Categories [Электронный ресурс].",
    "— URL: https : / / physbam .",
    "stanford .",
    "edu / cs448x / old / Procedural _ Noise(2f ) Categories.html (дата обращения: 11.11.2023).",
    "9.Крыжевич Л.",
    ",Ковалев В.",
    "Задача очистки изображения от шума и вейвлет-подходы к ее решению //Актуальные исследования в области математики, информатики, физики и методики их изучения в совре8Inменном образовательном пространстве.",
    "— 2016.",
    "1, № 1.",
    "39— 44.",
    "10.Смоленцев Н.",
    "Основы теории вейвлетов.",
    "Вейвлеты в MATLAB.",
    "— Москва : ДМК Пресс, 2014.",
    "628.",
    "11.Малла С.",
    "Вэйвлеты в обработке сигналов.",
    "— Москва : Мир, 2005.",
    "671.",
    "12.Новиков Л.",
    "Модифицированные вейвлеты в обработке данных аналитических приборов.",
    "I.",
    "основы теории.",
    "— 2006.",
    "— URL: https : / / cyberleninka .",
    "ru / article / n / modifitsirovannye - veyvlety - v obrabotke-dannyh-analiticheskih-priborov-i-osnovy-teorii (дата обращения: 25.11.2023).",
    "13.Meyer Y.",
    "Wavelets: Algorithms and applications.",
    "— Philadelphia : S.I.A.M., 1993.",
    "129.",
    "9In."
)

$n = $values.Length
for ($r = 1; $r -le $n; $r++) {
    $ws.Cells.Item($r, 1).Value = $values[$r - 1]
    if ($r -lt $n) {
        $ws.Cells.Item($r, 2).Value = $values[$r]
    } else {
        $ws.Cells.Item($r, 2).Value = $null
    }
}